$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the completed "1.- Redireccion de Logo a pagina principal..."
#    to-do item (the whole paragraph, including its paragraph mark).
# ---------------------------------------------------------------------------

# The paragraph carries the "_GoBack" bookmark reference at the very end of
# the to-do list (paragraph 13). Drop it first so it does not get silently
# reattached somewhere unexpected while we rearrange paragraphs.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$logoPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "1.- Redirecci*n de Logo a p*gina principal*") {
        $logoPara = $para
        break
    }
}
if ($logoPara -ne $null) {
    $logoPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Re-create the "_GoBack" bookmark at the start of what is now the first
#    numbered item ("2.- Front Quienes somos...").
# ---------------------------------------------------------------------------

$firstItemPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "2.- Front Quienes somos*") {
        $firstItemPara = $para
        break
    }
}
if ($firstItemPara -ne $null) {
    $bkStart = $firstItemPara.Range.Start
    $bkRange = $d.Range($bkStart, $bkStart)
    $d.Bookmarks.Add("_GoBack", $bkRange)
}

# ---------------------------------------------------------------------------
# 3) Item 13 ("Cambiar valores de Combobox en revisarMuestra"): mark
#    "revisarMuestra" as a proofing "spell start/end" span, same visible
#    text, just split into its own run.
# ---------------------------------------------------------------------------

$findRng = $d.Content
$findRng.Find.ClearFormatting()
if ($findRng.Find.Execute(" en revisarMuestra")) {
    $target = $d.Range($findRng.Start, $findRng.End)
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>revisarMuestra</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# 4) "Optional To Do:" heading: split "To" into its own proofing
#    spell-checked span; same visible text overall.
# ---------------------------------------------------------------------------

$findRng2 = $d.Content
$findRng2.Find.ClearFormatting()
if ($findRng2.Find.Execute(" To Do:")) {
    $target2 = $d.Range($findRng2.Start, $findRng2.End)
    $xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>To</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Do:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target2.InsertXML($xmlFrag2)
}
